$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format to avoid numeric auto-conversion/precision issues
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '67.107.27'
$ws.Range('E2').Value = '  +1.43%  '
$ws.Range('D3').Value = '3.856.65'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '468.47'
$ws.Range('E5').Value = '  +9.64%  '
$ws.Range('D6').Value = '144.83'
$ws.Range('E6').Value = '  +10.96%  '
$ws.Range('E7').Value = '  +3.41%  '
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').Value = '0.744'
$ws.Range('E9').Value = '  +2.23%  '
$ws.Range('D10').Value = '0.154'
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('D11').Value = '0.0000309'
$ws.Range('E11').Value = '  -8.00%  '
$ws.Range('D12').Value = '43.37'
$ws.Range('E12').Value = '  +4.31%  '
$ws.Range('D13').Value = '10.38'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('D14').Value = '4.480.58'
$ws.Range('E14').Value = '  +1.26%  '
$ws.Range('D15').Value = '14.82'
$ws.Range('E15').Value = '  -4.89%  '
$ws.Range('D16').Value = '3.863.94'
$ws.Range('E16').Value = '  +2.74%  '
$ws.Range('E17').Value = '  -0.37%  '
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('E19').Value = '  +6.09%  '
$ws.Range('D20').Value = '67.248.52'
$ws.Range('E20').Value = '  +1.41%  '
$ws.Range('D21').Value = '434.15'
$ws.Range('E21').Value = '  +4.56%  '
$ws.Range('D22').Value = '14.94'
$ws.Range('E22').Value = '  -1.11%  '
$ws.Range('D23').Value = '3.33'
$ws.Range('E23').Value = '  +6.79%  '
$ws.Range('D24').Value = '88.79'
$ws.Range('E24').Value = '  +4.55%  '
$ws.Range('E25').Value = '  +9.61%  '
$ws.Range('D26').Value = '37.93'
$ws.Range('E26').Value = '  +1.78%  '
$ws.Range('D27').Value = '10.11'
$ws.Range('E27').Value = '  +6.41%  '
$ws.Range('E28').Value = '  -2.33%  '
$ws.Range('D29').Value = '5.56'
$ws.Range('E29').Value = '  +2.93%  '
$ws.Range('D30').Value = '727.70'
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('D31').Value = '13.81'
$ws.Range('E31').Value = '  -0.58%  '
$ws.Range('E32').Value = '  +7.12%  '
$ws.Range('E33').Value = '  +0.75%  '
$ws.Range('D34').Value = '43.95'
$ws.Range('E34').Value = '  +13.18%  '
$ws.Range('E35').Value = '  +7.19%  '
$ws.Range('D36').Value = '57.97'
$ws.Range('E36').Value = '  +4.17%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').Value = '5.51'
$ws.Range('E38').Value = '  -3.21%  '
$ws.Range('D39').Value = '0.0483'
$ws.Range('E39').Value = '  +3.00%  '
$ws.Range('B40').Value = 'ThetaToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D40').Value = '2.93'
$ws.Range('E40').Value = '  +1.87%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').Value = '0.347'
$ws.Range('E41').Value = '  +8.21%  '
$ws.Range('E42').Value = '  +3.85%  '
$ws.Range('D43').Value = '0.0₃0677'
$ws.Range('E43').Value = '  -7.26%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').Value = '2.55'
$ws.Range('E45').Value = '  +5.95%  '
$ws.Range('D46').Value = '3.46'
$ws.Range('E46').Value = '  +2.05%  '
$ws.Range('D47').Value = '3.28'
$ws.Range('E47').Value = '  -0.56%  '
$ws.Range('D48').Value = '2.77'
$ws.Range('E48').Value = '  +5.21%  '
$ws.Range('E49').Value = '  +4.62%  '
$ws.Range('D50').Value = '2.89'
$ws.Range('E50').Value = '  +1.70%  '
$ws.Range('D51').Value = '144.10'
$ws.Range('E51').Value = '  +1.78%  '
